$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.908.28"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "1.845.79"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'309.27"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.4685"
$ws.Range("E7").Value = "  +3.42%  "
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "'0.07154"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "'0.9266"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("D11").Value = "'19.58"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "'0.07693"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.837.81"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'5.285"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "'6.407"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "'88.29"
$ws.Range("E16").Value = "  +3.57%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'0.000008638"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "26.951.03"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").Value = "'5.026"
$ws.Range("D23").Value = "'10.61"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "'1.923"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").Value = "'152.30"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").Value = "'2.015"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "'114.42"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").Value = "'4.877"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").Value = "'0.08846"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("E31").Value = "  +3.05%  "
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("D33").Value = "'0.7447"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").Value = "'2.794"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "'4.474"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").Value = "'0.01940"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'2.965"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").Value = "'0.05195"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").Value = "'0.5201"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("D41").Value = "'6.901"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").Value = "'0.1514"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").Value = "'8.140"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "'10.51"
$ws.Range("E44").Value = "  +5.07%  "
$ws.Range("D45").Value = "'0.4690"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "'1.006"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "'100.20"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'1.602"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").Value = "'65.62"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("D50").Value = "'0.06040"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'0.8920"
$ws.Range("E51").Value = "  +5.31%  "
